$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column B ("Estado civil") and column I ("Sexo") were reclassified from
# iaest-measure to iaest-dimension, so the metadata rows describing them
# need to be updated accordingly, and a new row of mapping-file references
# is appended.

# Row 3: semantic role URI (measure -> dimension)
$ws.Range("B3").Value = "iaest-dimension:estado-civil"
$ws.Range("I3").Value = "iaest-dimension:sexo"

# Row 4: kind (medida -> dim)
$ws.Range("B4").Value = "dim"
$ws.Range("I4").Value = "dim"

# Row 5: datatype (xsd:string -> skos:Concept)
$ws.Range("B5").Value = "skos:Concept"
$ws.Range("I5").Value = "skos:Concept"

# Row 6: new row with mapping file references for the two dimension columns
$ws.Range("B6").Value = "mapping-estado-civil.xlsx"
$ws.Range("I6").Value = "mapping-sexo.xlsx"

# Match the formatting used across the rest of the sheet by copying the
# cell format from an existing styled cell (same style as A1:I5)
$ws.Range("A5").Copy()
$ws.Range("B6").PasteSpecial(-4122)
$ws.Range("I6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wb.Save()
